# Add I0 and IF columns to the worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - values first
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-19
$iValues = @(7, 9, 5, 6, 7, 7, 7, 6, 5, 8, 9, 5, 7, 5, 2, 8, 7, 9)
$jValues = @(8, 9, 7, 7, 9, 8, 8, 7, 7, 8, 9, 5, 7, 5, 2, 8, 7, 9)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
